# Align PDF with App: sync data logic and remove legacy elements
#
# Recolors a set of "Output"/"Projections" feedback bullet paragraphs to the
# same red (EE0000) used throughout the rest of the document's review
# comments. Two of the target paragraphs already carried an explicit
# black/"text1" theme color; that is overwritten with the plain red RGB
# value (the theme reference is dropped), matching the rest.

$d = $word.ActiveDocument
$paras = $d.Paragraphs
$n = $paras.Count

# Distinctive substrings identifying each paragraph that must turn red.
$targets = @(
    "Add pie chart per entity",
    "Move this table to “Projections”",
    "Change “5. Projections” to “Projections” in label",
    "“Cash Flow Projections” that are under “Client & Structure” should be moved here",
    "Remove “Monte Carlo Wealth Projection”",
    "Replace “Select Model” with “Portfolio” and use drop down menu",
    "Remove “Analysis: Based on 1,000",
    "Is it possible to add a toggle “Before / After Tax”",
    "Is it possible to add a toggle “Nominal / Real”"
)

for ($i = 1; $i -le $n; $i++) {
    $p = $paras.Item($i)
    $t = $p.Range.Text
    foreach ($needle in $targets) {
        if ($t -like "*$needle*") {
            $p.Range.Font.Color = 238   # RGB(238,0,0) == EE0000
            break
        }
    }
}
